$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and two pairs of swapped rows)

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '20.713.93'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +2.14%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.513.69'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +4.27%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.21%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.9601'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.22%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '279.11'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.10%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3569'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.93%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3114'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.04%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.102'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +7.03%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.66'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.36%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.06720'
$ws.Range('D11').Style = "Normal"

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.9992'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.11%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '18.57'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +5.28%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.575'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +4.02%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.254'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.88%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.9619'
$ws.Range('D16').Style = "Normal"

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001027'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.97%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.505.44'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +4.10%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06044'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +5.83%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '70.00'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.35%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.579'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.33%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '14.80'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.31%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.31'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +4.71%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.283'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.14%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '20.779.57'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.25%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '146.33'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +3.92%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.147'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.87%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '17.47'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.91%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.668.22'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +4.28%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '115.90'
$ws.Range('D30').Style = "Normal"

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.986'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.17%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.047'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +4.77%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.8253'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +5.18%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.07992'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +3.30%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.203'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +7.34%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.455'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.47%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.05776'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.56%  '

$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.821'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +3.75%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.02058'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.56%  '

$ws.Range('E40').Value = '  +2.88%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9622'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.65%  '

$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1881'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.62%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '7.502'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.21%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.5298'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.24%  '

$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.536'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.89%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.29'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +4.54%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '120.50'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.26%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5283'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +3.49%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.856'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +6.94%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06487'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.59%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.9872'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.17%  '
